$d = $word.ActiveDocument

$replacements = @(
    @{Old = "314÷5="; New = "283÷2="},
    @{Old = "217÷7="; New = "576÷2="},
    @{Old = "969÷2="; New = "178÷2="},
    @{Old = "766÷2="; New = "447÷2="},
    @{Old = "875÷9="; New = "215÷3="},
    @{Old = "723÷4="; New = "909÷9="},
    @{Old = "975÷3="; New = "191÷8="},
    @{Old = "322÷4="; New = "483÷6="},
    @{Old = "468÷2="; New = "442÷4="},
    @{Old = "744÷7="; New = "423÷2="},
    @{Old = "497÷6="; New = "549÷2="},
    @{Old = "657÷3="; New = "104÷6="},
    @{Old = "935÷8="; New = "194÷7="},
    @{Old = "587÷2="; New = "270÷8="},
    @{Old = "783÷3="; New = "956÷2="},
    @{Old = "435÷3="; New = "593÷2="},
    @{Old = "115÷8="; New = "216÷7="},
    @{Old = "903÷8="; New = "262÷4="},
    @{Old = "607÷2="; New = "935÷6="},
    @{Old = "111÷2="; New = "539÷8="},
    @{Old = "160÷7="; New = "374÷3="},
    @{Old = "129÷4="; New = "293÷8="},
    @{Old = "416÷6="; New = "811÷3="},
    @{Old = "878÷2="; New = "183÷2="},
    @{Old = "651÷6="; New = "753÷5="}
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.New, 2)
}
